$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'isophonics_79'
$ws.Range("B2").Value = 'isophonics_297'
$ws.Range("C2").Value = 0.07435740514075886
$ws.Range("D2").Value = '[[''E'', ''B'', ''E'']]'
$ws.Range("E2").Value = '[[''G'', ''D'', ''G'']]'
$ws.Range("F2").Value = '[(4.103323, 13.890534)]'
$ws.Range("G2").Value = '[(0.421247, 3.083177)]'
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""

$ws.Range("A3").Value = 'schubert-winterreise_3'
$ws.Range("B3").Value = 'isophonics_33'
$ws.Range("C3").Value = 0.09821428571428571
$ws.Range("D3").Value = '[[''G#:maj'', ''C#:maj'', ''G#:maj'']]'
$ws.Range("E3").Value = '[[''G'', ''C'', ''G'']]'
$ws.Range("F3").Value = '[(9.64, 11.38)]'
$ws.Range("G3").Value = '[(13.723015, 16.39331)]'
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""

$ws.Range("A4").Value = 'schubert-winterreise_151'
$ws.Range("B4").Value = 'isophonics_156'
$ws.Range("C4").Value = 0.3882352941176471
$ws.Range("D4").Value = '[[''C:maj/G'', ''F:maj'', ''C:maj/G'']]'
$ws.Range("E4").Value = '[[''A'', ''D'', ''A'']]'
$ws.Range("F4").Value = '[(117.54, 121.8)]'
$ws.Range("G4").Value = '[(11.070127, 13.723731)]'
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""

$ws.Range("A5").Value = 'isophonics_133'
$ws.Range("B5").Value = 'isophonics_49'
$ws.Range("C5").Value = 0.2118055555555556
$ws.Range("D5").Value = '[[''E:min'', ''A'', ''D'', ''D/2'', ''D/3'']]'
$ws.Range("E5").Value = '[[''A:min'', ''D'', ''G'', ''G'', ''G'']]'
$ws.Range("F5").Value = '[(51.315598, 57.7109)]'
$ws.Range("G5").Value = '[(16.993365, 23.936132)]'
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""

$ws.Range("A6").Value = 'schubert-winterreise_106'
$ws.Range("B6").Value = 'schubert-winterreise_41'
$ws.Range("C6").Value = 0.2363636363636364
$ws.Range("D6").Value = '[[''B:min'', ''F#:7/A#'', ''B:min'']]'
$ws.Range("E6").Value = '[[''B:min'', ''F#:7/C#'', ''B:min/D'']]'
$ws.Range("F6").Value = '[(39.84, 47.28)]'
$ws.Range("G6").Value = '[(0.66, 2.68)]'
$ws.Range("H6").Value = 'spotify:track:1yerCi2iQCVkdHG6rdRn7R'
$ws.Range("I6").Value = ""

$ws.Range("A7").Value = 'isophonics_43'
$ws.Range("B7").Value = 'schubert-winterreise_191'
$ws.Range("C7").Value = 0.2870813397129187
$ws.Range("D7").Value = '[[''E'', ''A:min'', ''E'', ''A:min'']]'
$ws.Range("E7").Value = '[[''B:maj/F#'', ''E:min/G'', ''B:maj/F#'', ''E:min/G'']]'
$ws.Range("F7").Value = '[(31.051451, 41.036031)]'
$ws.Range("G7").Value = '[(87.24, 90.18)]'
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""

$ws.Range("A8").Value = 'schubert-winterreise_171'
$ws.Range("B8").Value = 'isophonics_180'
$ws.Range("C8").Value = 0.1001011122345804
$ws.Range("D8").Value = '[[''F#:maj'', ''C#:maj'', ''F#:maj'']]'
$ws.Range("E8").Value = '[[''F'', ''C'', ''F'']]'
$ws.Range("F8").Value = '[(46.56, 49.08)]'
$ws.Range("G8").Value = '[(0.440395, 5.558652)]'
$ws.Range("H8").Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'
$ws.Range("I8").Value = ""

$ws.Range("A9").Value = 'isophonics_136'
$ws.Range("B9").Value = 'isophonics_280'
$ws.Range("C9").Value = 0.1345050215208035
$ws.Range("D9").Value = '[[''E/4'', ''A'', ''A'', ''D/5'', ''A'']]'
$ws.Range("E9").Value = '[[''B'', ''E'', ''E/3'', ''A'', ''E/5'']]'
$ws.Range("F9").Value = '[(6.779, 11.273)]'
$ws.Range("G9").Value = '[(24.357755, 31.172811)]'
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""

$ws.Range("A10").Value = 'schubert-winterreise_114'
$ws.Range("B10").Value = 'jaah_43'
$ws.Range("C10").Value = 0.2232142857142857
$ws.Range("D10").Value = '[[''A:7/G'', ''D:maj/F#'', ''D:maj/A'']]'
$ws.Range("E10").Value = '[[''Bb:7'', ''Eb'', ''Eb'']]'
$ws.Range("F10").Value = '[(227.62, 235.06)]'
$ws.Range("G10").Value = '[(5.44, 10.01)]'
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = ""

$ws.Range("A11").Value = 'schubert-winterreise_160'
$ws.Range("B11").Value = 'schubert-winterreise_114'
$ws.Range("C11").Value = 0.2708333333333333
$ws.Range("D11").Value = '[[''D:maj'', ''G:maj'', ''D:maj/F#''], [''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range("E11").Value = '[[''D:maj/F#'', ''G:maj'', ''D:maj''], [''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range("F11").Value = '[(80.42, 84.64), (10.82, 17.08)]'
$ws.Range("G11").Value = '[(57.48, 64.58), (85.58, 94.24)]'
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""

$ws.Range("A12").Value = 'schubert-winterreise_61'
$ws.Range("B12").Value = 'schubert-winterreise_207'
$ws.Range("C12").Value = 0.3794871794871795
$ws.Range("D12").Value = '[[''D:7'', ''G:maj'', ''D:7/C'', ''G:maj/B'', ''D:7/C'', ''G:maj/B'']]'
$ws.Range("E12").Value = '[[''C:7'', ''F:maj'', ''C:7'', ''F:maj'', ''C:7'', ''F:maj'']]'
$ws.Range("F12").Value = '[(68.98, 87.2)]'
$ws.Range("G12").Value = '[(31.26, 45.88)]'
$ws.Range("H12").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I12").Value = ""

$ws.Range("A13").Value = 'schubert-winterreise_167'
$ws.Range("B13").Value = 'isophonics_194'
$ws.Range("C13").Value = 0.1030405405405405
$ws.Range("D13").Value = '[[''F:maj/G'', ''G:7'', ''C:maj'']]'
$ws.Range("E13").Value = '[[''G'', ''A:7'', ''D'']]'
$ws.Range("F13").Value = '[(27.66, 34.24)]'
$ws.Range("G13").Value = '[(73.125602, 78.292043)]'
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""

$ws.Range("A14").Value = 'schubert-winterreise_64'
$ws.Range("B14").Value = 'schubert-winterreise_65'
$ws.Range("C14").Value = 0.1348837209302326
$ws.Range("D14").Value = '[[''A:7'', ''D:min'', ''D:min''], [''D:maj/A'', ''A:7'', ''D:min''], [''D:min'', ''A:7/E'', ''D:min'']]'
$ws.Range("E14").Value = '[[''A#:7'', ''D#:min/A#'', ''D#:min''], [''D#:maj/A#'', ''A#:7'', ''D#:min/A#''], [''D#:min'', ''A#:7'', ''D#:min'']]'
$ws.Range("F14").Value = '[(12.42, 27.08), (55.72, 64.64), (0.82, 6.76)]'
$ws.Range("G14").Value = '[(36.18, 37.98), (35.88, 37.42), (9.1, 13.86)]'
$ws.Range("H14").Value = 'spotify:track:5UYEp9kllA47IhttiiMuJ0'
$ws.Range("I14").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

$ws.Range("A15").Value = 'isophonics_107'
$ws.Range("B15").Value = 'isophonics_31'
$ws.Range("C15").Value = 0.2289473684210526
$ws.Range("D15").Value = '[[''E'', ''A'', ''E'']]'
$ws.Range("E15").Value = '[[''G'', ''C'', ''G'']]'
$ws.Range("F15").Value = '[(15.027029, 21.737596)]'
$ws.Range("G15").Value = '[(17.581738, 30.271443)]'
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

$ws.Range("A16").Value = 'schubert-winterreise_37'
$ws.Range("B16").Value = 'schubert-winterreise_113'
$ws.Range("C16").Value = 0.1714285714285714
$ws.Range("D16").Value = '[[''F:min/C'', ''C'', ''F:min/C'']]'
$ws.Range("E16").Value = '[[''F:min'', ''C:maj'', ''F:min'']]'
$ws.Range("F16").Value = '[(45.58, 49.6)]'
$ws.Range("G16").Value = '[(0.78, 5.26)]'
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

$ws.Range("A17").Value = 'schubert-winterreise_53'
$ws.Range("B17").Value = 'jaah_27'
$ws.Range("C17").Value = 0.2053571428571428
$ws.Range("D17").Value = '[[''B:7/A'', ''E:maj/G#'', ''E:maj/B'']]'
$ws.Range("E17").Value = '[[''A:7'', ''D'', ''D'']]'
$ws.Range("F17").Value = '[(218.68, 226.54)]'
$ws.Range("G17").Value = '[(9.82, 12.93)]'
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
